$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Extend formatting from row 2 down through row 7 for columns A-F and H
#    (column G is never used in the data rows, so it is left alone).
# ---------------------------------------------------------------------------
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F7").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("H3:H7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Remove cells that must not exist in the final layout.
# ---------------------------------------------------------------------------
$ws.Range("A3").Clear()
$ws.Range("B4").Clear()
$ws.Range("C5").Clear()
$ws.Range("E6").Clear()

# ---------------------------------------------------------------------------
# 3. Populate the data, column by column (mirrors how the data was imported
#    so that the shared-strings table ends up ordered the same way).
# ---------------------------------------------------------------------------

# Column A - Request ID
$ws.Range("A2").Value = "HA-100520-HW2G5-PL-001"
$ws.Range("A4").Value = "GO-100520-HW3BT-PL-001"
$ws.Range("A5").Value = "HA-100520-HY4KJ-PL-008"
$ws.Range("A6").Value = "TL-100520-J23C3-PL-001"
$ws.Range("A7").Value = "HA-100520-HZ3ZJ-PL-004"

# Column B - Applicant Name
$ws.Range("B2").Value = "Tamanna . Vidhwani "
$ws.Range("B5").Value = "Girish . Naik "
$ws.Range("B6").Value = "Mayank Pratap. Tyagi "
$ws.Range("B7").Value = "Meher Deepthi . Rokkam "
$ws.Range("B3").Value = "TEJAS RAMESH, SHINGE "

# Column D - Address
$ws.Range("D2").Value = "09, Vidyut Nagar, Sindhi Colony- Gali No. 2, Patel Bhawan, Khandwa, 450001, Madhya Pradesh, India"
$ws.Range("D3").Value = "India, Maharashtra, Mumbai, 400092, 201, Khyati Deep CHS, Satya Sai Complex, Padma Nagar Road, Chikuwadi, Borivali (W)"
$ws.Range("D4").Value = "39/4, PWD COLONY, Jodhpur, 342001, India"
$ws.Range("D5").Value = "17,5A,SAHYADRI CO.HSG.SCO.,MUMBAI-PUNA ROAD, KALWA (W),THANE, Maharashtra, India-400605"
$ws.Range("D6").Value = "D-306, Shree Vardhman Mantra, Sector 67, Gurugram, 122102, Haryana, India"
$ws.Range("D7").Value = "Smt.Lalitha Shashtri Hostel,Rose Residency Complex,Sector-19,DWARKA, Delhi, India-110075"

# Column E - Fathers Name
$ws.Range("E2").Value = "Jayram Kundandas Vidhwani"
$ws.Range("E3").Value = "RAMESH TAYAPPA SHINGE"
$ws.Range("E4").Value = "Dinesh Kumar Mathur"
$ws.Range("E5").Value = "VASUDEV SITARAM NAIK"
$ws.Range("E7").Value = "Sudarsana Rao Rokkam"

# Column F - Customer type
$ws.Range("F2").Value = "General"
$ws.Range("F3").Value = "General"
$ws.Range("F4").Value = "General"
$ws.Range("F5").Value = "General"
$ws.Range("F6").Value = "General"
$ws.Range("F7").Value = "General"

# Column H - Country
$ws.Range("H2").Value = "India"
$ws.Range("H3").Value = "India"
$ws.Range("H4").Value = "India"
$ws.Range("H5").Value = "India"
$ws.Range("H6").Value = "India"
$ws.Range("H7").Value = "India"

# Column C - DOB (numeric dates, except row 7 which is free-text)
$ws.Range("C2").Value = 35597
$ws.Range("C3").Value = 31924
$ws.Range("C4").Value = 33723
$ws.Range("C6").Value = 34050
$ws.Range("C7").Value = "19-15-1988"

# ---------------------------------------------------------------------------
# 4. Column widths / best-fit sizing for columns A and D.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 26.6
$ws.Columns.Item(4).ColumnWidth = 119.0

# ---------------------------------------------------------------------------
# 5. Update the active selection to match the last-edited cell.
# ---------------------------------------------------------------------------
$ws.Range("C7").Select() | Out-Null
